# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (column G) holds a recomputed strike-count value for each
# saved game row. This updates the previously-written values with the
# freshly regenerated ones.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 3
    3  = 2
    4  = 3
    5  = 0
    6  = 1
    7  = 2
    8  = 0
    9  = 3
    10 = 0
    11 = 1
    12 = 2
    13 = 1
    14 = 1
    15 = 2
    16 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
